$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 320625
$ws.Range("D2").Value = 408622958
$ws.Range("C3").Value = 258
$ws.Range("D3").Value = 308479
$ws.Range("C4").Value = 316
$ws.Range("D4").Value = 452207
$ws.Range("C8").Value = 861
$ws.Range("D8").Value = 1266295
$ws.Range("C10").Value = 116962
$ws.Range("D10").Value = 171386841
$ws.Range("C12").Value = 59357
$ws.Range("D12").Value = 85669933
$ws.Range("C14").Value = 49
$ws.Range("D14").Value = 67543
$ws.Range("C16").Value = 4003
$ws.Range("D16").Value = 5681075
$ws.Range("C20").Value = 6643
$ws.Range("D20").Value = 9270246
$ws.Range("C22").Value = 77447
$ws.Range("D22").Value = 96577884
$ws.Range("C28").Value = 32453
$ws.Range("D28").Value = 47508425
$ws.Range("C30").Value = 11472
$ws.Range("D30").Value = 16500815
$ws.Range("C35").Value = 1826
$ws.Range("D35").Value = 2578410
$ws.Range("C36").Value = 97098
$ws.Range("D36").Value = 122215897
$ws.Range("C42").Value = 902
$ws.Range("D42").Value = 1327685
$ws.Range("C44").Value = 44356
$ws.Range("D44").Value = 65003278
$ws.Range("C46").Value = 9141
$ws.Range("D46").Value = 13118353
$ws.Range("C48").Value = 1404
$ws.Range("D48").Value = 1948603
$ws.Range("C51").Value = 2299
$ws.Range("D51").Value = 3208415
$ws.Range("C52").Value = 69041
$ws.Range("D52").Value = 86599294
$ws.Range("C57").Value = 382
$ws.Range("D57").Value = 561230
$ws.Range("C59").Value = 28181
$ws.Range("D59").Value = 41332157
$ws.Range("C62").Value = 11114
$ws.Range("D62").Value = 16071291
$ws.Range("C64").Value = 1358
$ws.Range("D64").Value = 1899289
$ws.Range("C68").Value = 1468
$ws.Range("D68").Value = 2056882
$ws.Range("C70").Value = 20494
$ws.Range("D70").Value = 26844470
$ws.Range("C73").Value = 60
$ws.Range("D73").Value = 87573
$ws.Range("C74").Value = 7580
$ws.Range("D74").Value = 11099008
$ws.Range("C76").Value = 5121
$ws.Range("D76").Value = 7434515
$ws.Range("C77").Value = 489
$ws.Range("D77").Value = 692239
$ws.Range("C78").Value = 276
$ws.Range("D78").Value = 387673
$ws.Range("C79").Value = 140830
$ws.Range("D79").Value = 175629798
$ws.Range("C83").Value = 428
$ws.Range("D83").Value = 624824
$ws.Range("C85").Value = 63497
$ws.Range("D85").Value = 93061232
$ws.Range("C86").Value = 81
$ws.Range("D86").Value = 120082
$ws.Range("C88").Value = 29706
$ws.Range("D88").Value = 42973422
$ws.Range("C90").Value = 2732
$ws.Range("D90").Value = 3933357
$ws.Range("C91").Value = 2824
$ws.Range("D91").Value = 3992864
$ws.Range("C92").Value = 33028
$ws.Range("D92").Value = 44748641
$ws.Range("C96").Value = 7980
$ws.Range("D96").Value = 11734288
$ws.Range("C98").Value = 7336
$ws.Range("D98").Value = 10640354
$ws.Range("C100").Value = 531
$ws.Range("D100").Value = 754716
$ws.Range("C101").Value = 494
$ws.Range("D101").Value = 712391
$ws.Range("C102").Value = 9614
$ws.Range("D102").Value = 14050818
$ws.Range("C104").Value = 2431
$ws.Range("D104").Value = 3803558
$ws.Range("C106").Value = 3256
$ws.Range("D106").Value = 5074796
$ws.Range("C108").Value = 143
$ws.Range("D108").Value = 223220
$ws.Range("C109").Value = 183
$ws.Range("D109").Value = 264543
$ws.Range("C110").Value = 141504
$ws.Range("D110").Value = 175002933
$ws.Range("C114").Value = 949
$ws.Range("D114").Value = 1393788
$ws.Range("C116").Value = 52677
$ws.Range("D116").Value = 77216646
$ws.Range("C118").Value = 27039
$ws.Range("D118").Value = 39173191
$ws.Range("C119").Value = 1311
$ws.Range("D119").Value = 1792784
$ws.Range("C122").Value = 2260
$ws.Range("D122").Value = 3174143
$ws.Range("C124").Value = 510118
$ws.Range("D124").Value = 673555727
$ws.Range("C125").Value = 91
$ws.Range("D125").Value = 120789
$ws.Range("C129").Value = 1375
$ws.Range("D129").Value = 2037986
$ws.Range("C131").Value = 207662
$ws.Range("D131").Value = 305269947
$ws.Range("C132").Value = 401
$ws.Range("D132").Value = 598250
$ws.Range("C134").Value = 182064
$ws.Range("D134").Value = 264716069
$ws.Range("C136").Value = 32
$ws.Range("D136").Value = 46832
$ws.Range("C137").Value = 2850
$ws.Range("D137").Value = 4004236
$ws.Range("C139").Value = 6338
$ws.Range("D139").Value = 8955307
$ws.Range("C142").Value = 44607
$ws.Range("D142").Value = 59553539
$ws.Range("C144").Value = 24
$ws.Range("D144").Value = 34730
$ws.Range("C148").Value = 14064
$ws.Range("D148").Value = 20625662
$ws.Range("C149").Value = 3755
$ws.Range("D149").Value = 5415751
$ws.Range("C152").Value = 401
$ws.Range("D152").Value = 576716
$ws.Range("C154").Value = 382
$ws.Range("D154").Value = 539663
$ws.Range("C155").Value = 17580
$ws.Range("D155").Value = 23237179
$ws.Range("C159").Value = 7162
$ws.Range("D159").Value = 10416985
$ws.Range("C161").Value = 4995
$ws.Range("D161").Value = 7188346
$ws.Range("C163").Value = 278
$ws.Range("D163").Value = 384431
$ws.Range("C164").Value = 266
$ws.Range("D164").Value = 380364
$ws.Range("C166").Value = 17110
$ws.Range("D166").Value = 26570937
$ws.Range("C167").Value = 1924
$ws.Range("D167").Value = 3049172
$ws.Range("C168").Value = 254
$ws.Range("D168").Value = 396802
$ws.Range("C170").Value = 58
$ws.Range("D170").Value = 94190
$ws.Range("C172").Value = 87715
$ws.Range("D172").Value = 109698302
$ws.Range("C177").Value = 641
$ws.Range("D177").Value = 944848
$ws.Range("C179").Value = 33840
$ws.Range("D179").Value = 49626554
$ws.Range("C181").Value = 13021
$ws.Range("D181").Value = 18815300
$ws.Range("C183").Value = 1245
$ws.Range("D183").Value = 1742429
$ws.Range("C185").Value = 1648
$ws.Range("D185").Value = 2316234
$ws.Range("C187").Value = 238563
$ws.Range("D187").Value = 296568415
$ws.Range("C193").Value = 877
$ws.Range("D193").Value = 1289845
$ws.Range("C195").Value = 86537
$ws.Range("D195").Value = 126854185
$ws.Range("C198").Value = 33026
$ws.Range("D198").Value = 47537691
$ws.Range("C200").Value = 16
$ws.Range("D200").Value = 22108
$ws.Range("C201").Value = 5109
$ws.Range("D201").Value = 7273022
$ws.Range("C204").Value = 4867
$ws.Range("D204").Value = 6740856
$ws.Range("C207").Value = 263853
$ws.Range("D207").Value = 326570180
$ws.Range("C209").Value = 254
$ws.Range("D209").Value = 363587
$ws.Range("C214").Value = 615
$ws.Range("D214").Value = 895378
$ws.Range("C216").Value = 95085
$ws.Range("D216").Value = 139106486
$ws.Range("C219").Value = 51421
$ws.Range("D219").Value = 74326568
$ws.Range("C220").Value = 32
$ws.Range("D220").Value = 45922
$ws.Range("C222").Value = 4672
$ws.Range("D222").Value = 6559305
$ws.Range("C225").Value = 5735
$ws.Range("D225").Value = 7933597
$ws.Range("C228").Value = 106305
$ws.Range("D228").Value = 132956279
$ws.Range("C230").Value = 74
$ws.Range("D230").Value = 106445
$ws.Range("C233").Value = 565
$ws.Range("D233").Value = 825439
$ws.Range("C235").Value = 49444
$ws.Range("D235").Value = 72431671
$ws.Range("C236").Value = 34
$ws.Range("D236").Value = 48711
$ws.Range("C237").Value = 12403
$ws.Range("D237").Value = 17836981
$ws.Range("C239").Value = 1894
$ws.Range("D239").Value = 2714882
$ws.Range("C241").Value = 2508
$ws.Range("D241").Value = 3508565
$ws.Range("C242").Value = 257642
$ws.Range("D242").Value = 325349607
$ws.Range("C244").Value = 250
$ws.Range("D244").Value = 359457
$ws.Range("C248").Value = 831
$ws.Range("D248").Value = 1220563
$ws.Range("C250").Value = 95746
$ws.Range("D250").Value = 140291814
$ws.Range("C253").Value = 65033
$ws.Range("D253").Value = 94245369
$ws.Range("C255").Value = 2412
$ws.Range("D255").Value = 3402009
$ws.Range("C258").Value = 4595
$ws.Range("D258").Value = 6452982
